$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.372.92"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.004.54"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "508.80"
$ws.Range("E5").Value = "  +2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.76"
$ws.Range("E6").Value = "  +3.86%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.429"
$ws.Range("E8").Value = "  +0.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.10"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.368"
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.516.10"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.30"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("E15").Value = "  +3.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.390.20"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.001.07"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.92"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.89"
$ws.Range("E19").Value = "  +3.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.00"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.71"
$ws.Range("E21").Value = "  +5.15%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.499"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.50"
$ws.Range("E24").Value = "  +2.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.124.21"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.165"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0936"
$ws.Range("E28").Value = "  +8.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.34"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.85"
$ws.Range("E30").Value = "  -1.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.79"
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.34"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.22"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("B35").Value = "EnergySwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.12"
$ws.Range("E35").Value = "  +14.22%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.44"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.81"
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0661"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.042.53"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.39"
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.77"
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.654"
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.197.77"
$ws.Range("E45").Value = "  +2.08%  "
$ws.Range("E46").Value = "  -0.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0240"
$ws.Range("E47").Value = "  +4.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.83"
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.920"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.58"
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0849"
$ws.Range("E51").Value = "  -0.73%  "
